$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blind 75")

# Row 19 - "Jump Game": mark completed, add solution note and runtime
$ws.Range("E19").Value = "X"
$ws.Range("F19").Value = "for every num if it can be reached up longest with max of longest and num+jump distance. If num cant be reached return false. Return true if last node reached."
$ws.Range("G19").Value = "O(N)"

# Row 22 - "Decode Ways": mark completed, add solution note and runtime
$ws.Range("E22").Value = "X"
$ws.Range("F22").Value = "create array of ans, for every num add the answer to the problem of num before. If 9<num and num before<27 add answer of 2 before as well. Return last in array of answers."
$ws.Range("G22").Value = "O(N)"

# Update the active sheet view/selection
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("I30").Select()
